# Fix image ratio issue in coco_plot_images:
# - reset the "contrast" default on test_settings from 1.1 to 1
# - give the new "min_size" column (H) an explicit "None" default on test_settings
# - make test_settings the active/selected sheet with H2 selected
$wb = $excel.ActiveWorkbook

$wsTest = $wb.Worksheets.Item("test_settings")

# contrast default: 1.1 -> 1
$wsTest.Range("A2").Value = 1

# min_size default for the new column: "None"
$wsTest.Range("H2").Value = "None"

# Make test_settings the active sheet and select H2 on it
$wsTest.Activate()
$wsTest.Range("H2").Select()
